# Keep straight quotes/apostrophes out of AutoFormat's way so replacement
# text matches the target's plain ASCII apostrophes.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false

$d = $word.ActiveDocument

# 1. Update the "Created" timestamp
$d.Content.Find.Execute(
    "Created April 16, 2025 at 16:26:59", $true, $false, $false, $false, $false,
    $true, 1, $false, "Created April 26, 2025 at 10:46:01", 2) | Out-Null

# 2. Remove the "The system overview is as follows:" paragraph and the blank
#    paragraph that follows it (they sit right after the "System Overview" heading).
$p = $d.Paragraphs(5)
$rng = $d.Range($p.Range.Start, $d.Paragraphs(7).Range.Start)
$rng.Delete() | Out-Null

# 3. Strip the leading "* " bullet markers and rename a couple of the
#    System Overview labels.
$d.Content.Find.Execute("* Date: 04-06-2025", $true, $false, $false, $false, $false,
    $true, 1, $false, "Date: 04-06-2025", 2) | Out-Null
$d.Content.Find.Execute("* Time: 16:02:09", $true, $false, $false, $false, $false,
    $true, 1, $false, "Time: 16:02:09", 2) | Out-Null
$d.Content.Find.Execute("* Operating System Name: kb322-18", $true, $false, $false, $false, $false,
    $true, 1, $false, "OS Name: kb322-18", 2) | Out-Null
$d.Content.Find.Execute(
    "* Operating System Version: #1 SMP PREEMPT_DYNAMIC Debian 6.1.129-1 (2025-03-06)",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "OS Version: #1 SMP PREEMPT_DYNAMIC Debian 6.1.129-1 (2025-03-06)", 2) | Out-Null
$d.Content.Find.Execute("* Computer Name: kb322-18", $true, $false, $false, $false, $false,
    $true, 1, $false, "Computer Name: kb322-18", 2) | Out-Null
$d.Content.Find.Execute("* IP Address: 140.160.138.147", $true, $false, $false, $false, $false,
    $true, 1, $false, "IP Address: 140.160.138.147", 2) | Out-Null

# 4. Patch Status Summary: collapse the pending-patch paragraphs into a
#    single "no updates" line, keeping the trailing blank paragraph intact.
$startPara = $null
$endPara = $null
$i = 0
foreach ($para in $d.Paragraphs) {
    $i = $i + 1
    $t = $para.Range.Text
    if ($t -eq "The following patch is pending update:`r") { $startPara = $i }
    if ($t -eq "A vulnerability exists in an unknown function of a file within one of the affected products, which could lead to remote attack and SQL injection.`r") { $endPara = $i }
}
$start = $d.Paragraphs($startPara).Range.Start
$end = $d.Paragraphs($endPara + 1).Range.Start
$rng = $d.Range($start, $end)
$rng.Delete() | Out-Null
$rngIns = $d.Range($start, $start)
$rngIns.InsertParagraphAfter() | Out-Null
$rngTxt = $d.Range($start, $start)
$rngTxt.Text = "There are no pending updates available."

# 5. Compliance with RMF Controls: collapse the remediation bullets into a
#    single paragraph, keeping the trailing blank paragraph intact.
$startPara = $null
$endPara = $null
$i = 0
foreach ($para in $d.Paragraphs) {
    $i = $i + 1
    $t = $para.Range.Text
    if ($t -eq "For remediation of the identified vulnerability:`r") { $startPara = $i }
    if ($t -eq "* Vulnerability Checks: Regularly run vulnerability scans to detect any potential security breaches.`r") { $endPara = $i }
}
$start = $d.Paragraphs($startPara).Range.Start
$end = $d.Paragraphs($endPara + 1).Range.Start
$rng = $d.Range($start, $end)
$rng.Delete() | Out-Null
$rngIns = $d.Range($start, $start)
$rngIns.InsertParagraphAfter() | Out-Null
$rngTxt = $d.Range($start, $start)
$rngTxt.Text = "To ensure compliance, it is essential to have a process in place for identifying and reporting vulnerabilities. The current system does not have any identified patches that require remediation."

# 6. Recommended next steps: rewrite the three numbered steps.
$d.Content.Find.Execute(
    "1. Review and Assess Updates: Carefully review the pending update to understand its impact on the system's security.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Review the current patch status and schedule any necessary patch deployments.", 2) | Out-Null
$d.Content.Find.Execute(
    "2. Scheduling patch deployments, if needed: If the review confirms that the update is necessary, schedule its deployment at a suitable time when the system can be isolated from potential threats.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Update documentation should be reviewed and updated to reflect the current patch status.", 2) | Out-Null
$d.Content.Find.Execute(
    "3. Guidance for Update documentation: Ensure that detailed records are maintained of all updates made to the system, including dates, times, and personnel involved.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Schedule regular review and assessment of updates to ensure the system remains secure.", 2) | Out-Null

# 7. Risk Assessment paragraph rewrite.
$d.Content.Find.Execute(
    "The pending patch contains an unknown function within one of the affected products, which could lead to remote attack and SQL injection, potentially causing disruption to the system or revealing sensitive information. To mitigate this risk, it is recommended that the update be applied as soon as possible to prevent potential exploitation.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Since there are no pending updates, the risk level is considered low. There is no potential impact on the system's security at this time. However, it is still essential to regularly review and assess updates to ensure the system remains secure.", 2) | Out-Null
